$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 95; this shifts existing rows 95..168 down to 96..169
$ws.Rows.Item(95).Insert()

# Populate the new row 95 with the new data point
$ws.Cells.Item(95, 1).Value = 10
$ws.Cells.Item(95, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(95, 3).Value = "La Araucanía"
$ws.Cells.Item(95, 4).Value = 44574
$ws.Cells.Item(95, 5).Value = 9
$ws.Cells.Item(95, 6).Value = 100112005
$ws.Cells.Item(95, 7).Value = "Puerro"
$ws.Cells.Item(95, 8).Value = "Azul de Maquehue"
$ws.Cells.Item(95, 9).Value = "Primera"
$ws.Cells.Item(95, 10).Value = 55
$ws.Cells.Item(95, 11).Value = 14000
$ws.Cells.Item(95, 12).Value = 14000
$ws.Cells.Item(95, 13).Value = 14000
$ws.Cells.Item(95, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(95, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(95, 16).Value = 1167
$ws.Cells.Item(95, 17).Value = 12
$ws.Cells.Item(95, 18).Value = "Hortaliza"
